$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.362.19'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -4.07%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.862.34'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -4.94%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -1.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '323.86'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.34%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.98%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4524'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -5.77%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3876'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -5.27%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '48.17'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -10.87%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07913'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.023'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -3.56%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '21.47'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -4.61%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.856.15'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -5.20%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.906'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.158'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -5.74%  '
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.22%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001032'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -3.58%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '85.76'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -5.48%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06508'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.74%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.15'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -7.45%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.001'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.97%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.532'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -5.41%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.352.28'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -4.12%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.88'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -5.19%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.275'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -1.15%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.068.59'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -5.49%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '152.92'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -2.43%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.81'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -2.70%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.063'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -5.37%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -5.98%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '120.93'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -2.86%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.499'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +2.99%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09316'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -3.79%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9362'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -5.41%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -2.32%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.277'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -6.41%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02239'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -4.26%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.38%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06002'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -3.28%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -9.92%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9999'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -1.00%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5911'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -5.36%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1897'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.36%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.15'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -9.59%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.274'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -5.26%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5647'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -5.36%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '12.05'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -7.20%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.372'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.17%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.927'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -6.58%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06787'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.47%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '108.11'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -3.00%  '
